# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (A1) onto the new
# header cells so they match the bold/centered/bordered look of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values, repeated for every player row (2-47).
$ws.Range("AD2:AD47").Value = 79
$ws.Range("AE2:AE47").Value = 83
$ws.Range("AF2:AF47").Value = 0
